$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct district names to official names (per commit message)
$ws.Range("G3").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G4").Value = "Davangere"
$ws.Range("G6").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G7").Value = "Vijayapura (Bijapur)"
$ws.Range("G9").Value = "Davangere"
$ws.Range("G18").Value = "Davangere"
$ws.Range("G20").Value = "Vijayapura (Bijapur)"
$ws.Range("G22").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G29").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G33").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G35").Value = "Bagalkot"
$ws.Range("G38").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G42").Value = "Vijayapura (Bijapur)"
$ws.Range("G47").Value = "Davangere"
$ws.Range("G50").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G51").Value = "Bagalkot"
$ws.Range("G52").Value = "Davangere"
$ws.Range("G54").Value = "Kalaburagi (Gulbarga)"

# Remove stray empty inline-string cells in column F
$ws.Range("F10").ClearContents()
$ws.Range("F21").ClearContents()
$ws.Range("F31").ClearContents()
